$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellref, $val) {
    $r = $ws.Range($cellref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '23.491.20'
$ws.Range('E2').Value = '  -0.12%  '

$ws.Range('D3').Value = '1.650.64'
$ws.Range('E3').Value = '  -0.23%  '

Set-CellText 'D4' '1.000'
$ws.Range('E4').Value = '  -0.12%  '

Set-CellText 'D5' '1.000'
$ws.Range('E5').Value = '  -0.06%  '

Set-CellText 'D6' '300.39'
$ws.Range('E6').Value = '  -0.63%  '

Set-CellText 'D7' '0.3786'
$ws.Range('E7').Value = '  -1.28%  '

Set-CellText 'D8' '50.62'
$ws.Range('E8').Value = '  -1.07%  '

Set-CellText 'D9' '0.3502'
$ws.Range('E9').Value = '  -2.48%  '

Set-CellText 'D10' '1.226'
$ws.Range('E10').Value = '  -1.24%  '

Set-CellText 'D11' '0.08059'
$ws.Range('E11').Value = '  -1.58%  '

Set-CellText 'D12' '1.000'
$ws.Range('E12').Value = '  -0.12%  '

Set-CellText 'D13' '22.11'
$ws.Range('E13').Value = '  -1.39%  '

Set-CellText 'D14' '6.321'
$ws.Range('E14').Value = '  -2.64%  '

Set-CellText 'D15' '7.271'
$ws.Range('E15').Value = '  -2.94%  '

Set-CellText 'D16' '0.00001212'
$ws.Range('E16').Value = '  -0.63%  '

$ws.Range('D17').Value = '1.650.92'
$ws.Range('E17').Value = '  -0.14%  '

Set-CellText 'D18' '95.44'
$ws.Range('E18').Value = '  -2.09%  '

Set-CellText 'D19' '0.06973'
$ws.Range('E19').Value = '  -0.04%  '

Set-CellText 'D20' '6.639'
$ws.Range('E20').Value = '  -2.65%  '

Set-CellText 'D21' '17.46'
$ws.Range('E21').Value = '  -1.06%  '

Set-CellText 'D22' '1.000'
$ws.Range('E22').Value = '  -0.11%  '

Set-CellText 'D23' '12.48'
$ws.Range('E23').Value = '  -1.55%  '

$ws.Range('D24').Value = '23.490.04'
$ws.Range('E24').Value = '  -0.17%  '

Set-CellText 'D25' '2.425'
$ws.Range('E25').Value = '  -3.10%  '

Set-CellText 'D26' '3.024'
$ws.Range('E26').Value = '  +0.75%  '

Set-CellText 'D27' '21.11'
$ws.Range('E27').Value = '  -0.55%  '

Set-CellText 'D28' '151.78'
$ws.Range('E28').Value = '  -0.15%  '

Set-CellText 'D29' '5.174'
$ws.Range('E29').Value = '  -1.29%  '

Set-CellText 'D30' '131.85'
$ws.Range('E30').Value = '  -1.49%  '

$ws.Range('D31').Value = '1.837.37'
$ws.Range('E31').Value = '  +0.00%  '

Set-CellText 'D32' '6.910'
$ws.Range('E32').Value = '  -4.03%  '

Set-CellText 'D33' '2.140'
$ws.Range('E33').Value = '  -4.87%  '

Set-CellText 'D34' '11.21'
$ws.Range('E34').Value = '  -7.27%  '

Set-CellText 'D35' '0.9902'
$ws.Range('E35').Value = '  -6.31%  '

$ws.Range('E36').Value = '  -2.88%  '

Set-CellText 'D37' '0.08769'
$ws.Range('E37').Value = '  -0.25%  '

Set-CellText 'D38' '5.946'
$ws.Range('E38').Value = '  -3.01%  '

$ws.Range('E39').Value = '  -2.88%  '

Set-CellText 'D40' '0.06838'
$ws.Range('E40').Value = '  -2.63%  '

Set-CellText 'D41' '12.91'
$ws.Range('E41').Value = '  -2.35%  '

Set-CellText 'D42' '0.6914'
$ws.Range('E42').Value = '  -1.43%  '

Set-CellText 'D43' '1.296'
$ws.Range('E43').Value = '  -2.97%  '

Set-CellText 'D44' '15.52'
$ws.Range('E44').Value = '  -2.94%  '

$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-CellText 'D45' '0.9993'
$ws.Range('E45').Value = '  -0.12%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-CellText 'D46' '0.6401'
$ws.Range('E46').Value = '  -2.00%  '

Set-CellText 'D47' '2.254'
$ws.Range('E47').Value = '  -2.33%  '

$ws.Range('E48').Value = '  -0.86%  '

Set-CellText 'D49' '127.20'
$ws.Range('E49').Value = '  -0.51%  '

Set-CellText 'D50' '0.07683'
$ws.Range('E50').Value = '  -2.93%  '

Set-CellText 'D51' '1.241'
$ws.Range('E51').Value = '  +3.21%  '
